$d = $word.ActiveDocument

# --- Edit 1 --------------------------------------------------------------
# Merge the runs around "roof;" (removing the gramStart/gramEnd proofErr
# markers) so the sentence becomes a single run.
$old1 = "The rain was still coming down as a dull roar on the roof; blending with the static from the television. "
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2)

# --- Edit 2 ----------------------------------------------------------------
# Merge the runs around "couldn't" in the "Looking at the clock on the wall"
# paragraph into a single run (removing the proofErr markers).
$old2 = "Looking at the clock on the wall I realized it was suddenly 8pm! I decided to make some chamomile tea and go to bed after. I couldn’t enjoy the weekend off if I kept falling asleep and sleeping in my bed was much better on my back than the recliner. My phone should have enough battery to use the light to read a little more until I doze off again. I stumbled back into the kitchen, still half asleep and began boiling some water."
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2)

# --- Edit 3 ----------------------------------------------------------------
# Merge the runs around "didn't" in the "Glancing outside" paragraph into a
# single run (removing the proofErr markers).
$old3 = " since; I had been out for several hours at this point. I shivered a little and drew the curtains shut as the kettle began to boil. Crash! I didn’t see the lightning, but the thunderclap sounded as if it was right above the house. A few minutes later I returned to the living room, steaming mug in hand. "
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)

# --- Edit 4 ----------------------------------------------------------------
# Merge the runs around "couldn't" in the "they were just gone" paragraph
# into a single run, then fix the "dak" -> "dark" typo while keeping the
# inserted "r" in its own run (split off from its neighbours).
$old4 = "I couldn’t tell if it was simply too dak to see that far or if my house had been teleported to some other dimension that was pitch black and full of nothing but swirling storm clouds, thunder, and rain. I suppose in the end it "
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2)

$rng = $d.Content
$found4b = $rng.Find.Execute("too dak", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $rng.Start + 6

$insertRange = $d.Range($pos, $pos)
$insertRange.InsertAfter("r")
$rSpan = $d.Range($pos, $pos + 1)
$rSpan.Bold = 1
$rSpan.Bold = 0

Write-Host "Edit1:" $found1 "Edit2:" $found2 "Edit3:" $found3 "Edit4:" $found4 "Edit4b:" $found4b
